$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in A2:B53 (Date / Ratio) is currently sorted by Date descending.
# Re-sort ascending by Date (column A), matching Data > Sort in Excel.
$sortRange = $ws.Range("A2:B53")
$keyRange  = $ws.Range("A2:A53")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Selection afterwards covers the full A:B columns.
$ws.Range("A1:B1048576").Select()
